$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update prices in column D for rows 28-33
$ws.Range("D28").Value = 7530
$ws.Range("D29").Value = 7950
$ws.Range("D30").Value = 9500
$ws.Range("D31").Value = 10420
$ws.Range("D32").Value = 10830
$ws.Range("D33").Value = 11550
